$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "protein_meal"
$ws.Range("B11").Value = "Protein Meal"

$ws.Range("B11").Select()
